# tambah table master agent
$wb = $excel.ActiveWorkbook

# --- Master sheet: add new row "master agent" / "MST010" ---
$wsMaster = $wb.Worksheets.Item("Master")
$wsMaster.Range("A7").Value = "master agent"
$wsMaster.Range("B7").Value = "MST010"
$null = $wsMaster.Range("A8").Select()

# --- Transaksi sheet: swap A2/B2 values, adjust column widths, move selection ---
$wsTransaksi = $wb.Worksheets.Item("Transaksi")
$valA2 = $wsTransaksi.Range("A2").Value2
$valB2 = $wsTransaksi.Range("B2").Value2
$wsTransaksi.Range("A2").Value = $valB2
$wsTransaksi.Range("B2").Value = $valA2
$wsTransaksi.Columns.Item(1).ColumnWidth = 21.6
$null = $wsTransaksi.Range("B14").Select()
